$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Proofreading pass: re-typed cell text in columns A and C (and two cells
# in column F) to drop embedded line-breaks / tidy wording. Order matches the
# author's edit sequence: column A top-to-bottom, then column C, then F.
$ws.Range("A2").Value2  = "Visual Light Sensitivity  Questionnaire-8"
$ws.Range("A9").Value2  = "The Pittsburgh Sleep Quality  Index (PSQI)"
$ws.Range("A10").Value2 = "Self-Rating of Biological Rhythm Disorder for Disorder for  Adolescents (SBRDA)"

$ws.Range("C2").Value2  = "Eight-question survey  to assess the  presence and severity of photosensitivity symptoms"
$ws.Range("C3").Value2  = "30 items survey  to assess electrical lighting environment in office"
$ws.Range("C5").Value2  = "23 items questionnaire to assess light environment  in a hospital"
$ws.Range("C7").Value2  = "17 items questionnaire to understand individuals  phase of  entrainment"
$ws.Range("C8").Value2  = "16 Factor questionnaire measuring  practice, behavior and attitude  related sleep"
$ws.Range("C9").Value2  = "9 items inventory  to measure sleep  quality and sleeping pattern"
$ws.Range("C10").Value2 = "29 Items questionnaire  assessing four  dimensions of biological rhythm  disorder in adolescents "
$ws.Range("C11").Value2 = "16 dichotomous items  questionnaire to assess  ""photophobia"" and ""photophilia"""

$ws.Range("F6").Value2  = "Correlation with oral temperature"
$ws.Range("F9").Value2  = "Correlation  with  clinical measurements"

# --- Column A is no longer auto-fit to its longest entry; set an explicit
# manual width instead (Excel's ColumnWidth differs from the stored sheet
# width by the default padding, ~0.834 chars).
$ws.Columns.Item(1).ColumnWidth = 43.6659375

# --- View state: zoomed to 110% with D10 as the active/selected cell.
[void]$ws.Range("D10").Select()
$excel.ActiveWindow.Zoom = 110
